$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "joinable_relation.id"
$ws.Range("E1").Value = "joinable_relation.foreign_field"
$ws.Range("F1").Value = "joinable_relation.another_foreign_field"
$ws.Range("G1").Value = "another_joinable_relation.id"
$ws.Range("H1").Value = "another_joinable_relation.foreign_field"
$ws.Range("I1").Value = "another_joinable_relation.another_foreign_field"

$ws.Range("I1").Select()
